$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = -0.02593244657879423
$ws.Range("B3").Value = 0.8322916158076636
$ws.Range("B4").Value = 65.00800463600758
$ws.Range("B5").Value = "[1.0, 0.25258696009044634, 0.21797593258860617, 0.2332947276494444, 0.22263181012358957, 0.4168504398817456, 0.23623188868210593, 0.2143796725892713, 0.19424242364604458, 0.19867236992249107, 0.15698687759633312, 0.17488772338702271, 0.15305503121444938, 0.17865630752913067, 0.20541856834833336, 0.16363777366356688, 0.14998147432887562, 0.14611546778720555, 0.12397716902668973, 0.15079325822268036]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -2.32256187444717
$ws.Range("B5").Value = 0.7840915543792709
$ws.Range("B6").Value = 53.61790160812556
$ws.Range("B7").Value = 2.499771463899311
$ws.Range("B8").Value = 0.7286996154842171
$ws.Range("B9").Value = 41.49333159035533
$ws.Range("B10").Value = "[1.0, 0.2590163508111495, 0.21849322228450127, 0.22919936749035894, 0.26149467565769147, 0.3237745069592172, 0.2185838817942732, 0.19955513693769017, 0.17664824043633395, 0.2086154968372668, 0.16515991718540746, 0.16868113179261457, 0.16232640224493666, 0.15967328714459456, 0.1622265011295381, 0.13378248466072254, 0.14120832897673705, 0.15402239310562757, 0.13418677445071722, 0.14689705531848932]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.01959257089766838
$ws.Range("B3").Value = 0.2437428974517225
$ws.Range("B4").Value = 0.08196289620046414
$ws.Range("B5").Value = 0.9180371518230247
$ws.Range("B6").Value = "[1.0, -0.0002703759838989657, -0.021003629998874195, -0.032617266914086535, -0.038620045378157525, 0.3104179679214549, -0.0021862482309711745, 0.024676328814873354, -0.0009056348000959117, -0.0038641263678660747, -0.06594555141554542, -0.04778972823700245, -0.038568537814385885, -0.002503237336769498, 0.039306308591957405, -0.017745604589805102, -0.023602582530362272, -0.04465798964643314, -0.060816197545603266, -0.028669702033884727]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.02404507333627934
$ws.Range("B3").Value = 0.247176943251281
$ws.Range("B4").Value = 0.06636574118345585
$ws.Range("B5").Value = 0.9190637193571399
$ws.Range("B6").Value = "[0.9999999999999999, 0.0004676088599647633, -0.02095352733807903, -0.03316525682029443, -0.038599037969873524, 0.3099755590604995, -0.0011973071399435974, 0.024839079563204036, -0.0004152391501257935, -0.003044124289353733, -0.06608402925502255, -0.047424735261551554, -0.037995584649338096, -0.0019149817943322522, 0.04066764378264819, -0.01707294139507888, -0.023303672927699078, -0.0440603652206585, -0.060565755821467206, -0.0279868823823532]"
$ws.Range("B7").Value = 0.02914107880281368

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.04540099676221379
$ws.Range("B3").Value = 0.2035823066711887
$ws.Range("B4").Value = 0.07595619762525964
$ws.Range("B5").Value = 0.9241694488609393
$ws.Range("B6").Value = "[0.9999999999999999, 0.020218052274143583, -0.017987560166422178, 0.001320742084632113, -0.01593790457703244, 0.24572072193707523, -0.0069903574365254165, 0.0191817852163228, -0.02346873290582793, -0.003944662486433537, -0.04763049205054739, -0.03430387037730003, -0.03567737726458859, -0.013160551237106206, 0.04516696492157638, -0.027635241900779768, -0.01679397124900664, -0.026895261925194355, -0.05064239841025293, -0.018577964541152073]"
$ws.Range("B7").Value = -0.0002512753242768607
$ws.Range("B9").Value = 0.8178229721536755
